# The source workbook has a single worksheet named "merged" (produced by a
# pandas/openpyxl export). Rename it back to Excel's default "Sheet1", which
# is the only cell/sheet-level change between the before/after workbooks.
$wb = $excel.ActiveWorkbook

$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "merged") {
        $ws = $sheet
        break
    }
}
if ($ws -eq $null) {
    $ws = $wb.ActiveSheet
}

$ws.Name = "Sheet1"
